# LocalPublicResults.pptx update
#
# - Move/resize "TextBox 137" (owlcms / IP address label) slightly.
# - Change the laptop's displayed IP from 192.168.1.116 to 192.168.4.100.
# - Shrink the curved connector that points at that textbox (154) to match
#   the new textbox position.
# - Update the "Competition Network (192.168.1.x)" caption to
#   "Competition Network (192.168.4.x)".
# - Shrink the other curved connector (61) that points at the same textbox.
#
# NOTE: PowerPoint's Shape.Left/Top/Width/Height are expressed in points
# (1 pt = 12700 EMU) and the host truncates (floors) the point value times
# 12700 when it re-serialises to EMU, so literal point values are chosen to
# land in the middle of the EMU's valid point interval to avoid rounding
# down to the EMU below the intended target.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) "TextBox 137" (owlcms / 192.168.x.x label): reposition + widen by 1 EMU.
$tbIp = $s.Shapes.Item("TextBox 137")
$tbIp.Left = 340.19122047244093   # 4320428 EMU (was 4397375 EMU / 346.25 pt)
$tbIp.Width = 85.47681102362205   # 1085555 EMU (was 1085554 EMU)

# 2) Update the IP address text run inside that textbox, leaving the
#    "owlcms" run and the line break untouched.
$tr = $tbIp.TextFrame.TextRange
$ipLen = "192.168.1.116".Length
$ipStart = $tr.Length - $ipLen + 1
$tr.Characters($ipStart, $ipLen).Text = "192.168.4.100"

# 3) Curved connector "Connector: Curved 154" shrinks horizontally to track
#    the textbox's new left edge.
$conn154 = $s.Shapes.Item("Connector: Curved 154")
$conn154.Width = 211.42964566929135   # 2685156 EMU (was 2762102 EMU)

# 4) "TextBox 160": update the second run, " Network (192.168.1.x)" ->
#    " Network (192.168.4.x)"; the first run ("Competition") is untouched.
$tbNet = $s.Shapes.Item("TextBox 160")
$trNet = $tbNet.TextFrame.TextRange
$oldSuffix = " Network (192.168.1.x)"
$suffixLen = $oldSuffix.Length
$suffixStart = $trNet.Length - $suffixLen + 1
$trNet.Characters($suffixStart, $suffixLen).Text = " Network (192.168.4.x)"

# 5) Curved connector "Connector: Curved 61" also shrinks horizontally.
$conn61 = $s.Shapes.Item("Connector: Curved 61")
$conn61.Width = 52.17964566929134   # 662681 EMU (was 739627 EMU)
